$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy the date format from A58 onto the new date cells A59:A61
$ws.Range("A58").Copy() | Out-Null
$ws.Range("A59:A61").PasteSpecial(-4122) | Out-Null  # xlPasteFormats

# Insert a new row of data (2026-01-27, 48) before the existing last row, then shift the
# old 46048/66 row down, and append a new trailing row (2026-01-28, 117).
$ws.Range("A59").Value = 46049
$ws.Range("B59").Value = 48

$ws.Range("A60").Value = 46048
$ws.Range("B60").Value = 66

$ws.Range("A61").Value = 46050
$ws.Range("B61").Value = 117

# Update the active selection to match the author's final cursor position.
$ws.Range("F59").Select() | Out-Null
